$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of data describing the LCD's max clock rate
$ws.Range("A8").Value = "Max clock for lcd"
$ws.Range("B8").Value = "25.641 khz"
$ws.Range("B8").HorizontalAlignment = -4108  # xlCenter

$ws.Range("B15").Select()
